$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell H1 = "Save", matching the formatting of the other header cells (e.g. G1)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# "Save" values for data rows 2..12
$saveValues = @(0, 0, 0, 0, 0, 0, 0, 1, 0, 1, 1)
for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
